$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column Q (year 2023) with the same look/formatting as column P (2022) ---
# Copy formats from column P (rows 3..25) into the new column Q first, then fill in values.
$ws.Range("P3:P25").Copy() | Out-Null
$ws.Range("Q3:Q25").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
# Row 4 never had a P cell to begin with (no 2022 value there), so undo the stray paste.
$ws.Range("Q4").Clear() | Out-Null

# New data values for 2023 (column Q)
$ws.Range("Q3").Value = 2023

$ws.Range("Q5").Value = 2385.9
$ws.Range("Q6").Value = 112.1
$ws.Range("Q7").Value = 267.89999999999998
$ws.Range("Q8").Value = 230.9
$ws.Range("Q9").Value = 249.7
$ws.Range("Q10").Value = 287
$ws.Range("Q11").Value = 334.7
$ws.Range("Q12").Value = 851
$ws.Range("Q13").Value = 48.5
$ws.Range("Q14").Value = 4.2
# Row 15 stays blank (format only, already copied above)

$ws.Range("Q16").Value = 26.890545708088244
$ws.Range("Q17").Value = 15.490056759274875
$ws.Range("Q18").Value = 22.218388220841799
$ws.Range("Q19").Value = 29.614327895683314
$ws.Range("Q20").Value = 30.104452089276922
$ws.Range("Q21").Value = 21.825966598728439
$ws.Range("Q22").Value = 32.351574864874735
$ws.Range("Q23").Value = 30.810022297218843
$ws.Range("Q24").Value = 29.193884213235311
$ws.Range("Q25").Value = 7.4362892319581295

# --- Column A:C width tweak (37.140625 -> 36.28515625 serialized units) ---
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 35.5

# --- Rows 4-25 get an explicit (custom) row height of 15 ---
$ws.Rows("4:25").RowHeight = 15

# --- Reset the saved selection back to the default cell ---
$ws.Range("A1").Select() | Out-Null
